$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 388-389, pushing the existing rows 388-393
# down to 390-395 (same row content, unchanged).
$ws.Rows("388:389").Insert()

# Row 388 (new): 1a plateado, 3/Feb/2022, Provincia de Melipilla
$ws.Cells.Item(388, 1).Value = 4
$ws.Cells.Item(388, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(388, 3).Value = "Los Lagos"
$ws.Cells.Item(388, 4).Value = 44595
$ws.Cells.Item(388, 5).Value = 10
$ws.Cells.Item(388, 6).Value = "Fruta"
$ws.Cells.Item(388, 7).Value = 100102
$ws.Cells.Item(388, 8).Value = "Cítricos"
$ws.Cells.Item(388, 9).Value = 100102003
$ws.Cells.Item(388, 10).Value = "Limón"
$ws.Cells.Item(388, 11).Value = "Sin especificar"
$ws.Cells.Item(388, 12).Value = "1a plateado"
$ws.Cells.Item(388, 13).Value = 600
$ws.Cells.Item(388, 14).Value = 23000
$ws.Cells.Item(388, 15).Value = 24000
$ws.Cells.Item(388, 16).Value = 23500
$ws.Cells.Item(388, 17).Value = "$/malla 18 kilos"
$ws.Cells.Item(388, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(388, 19).Value = 1306
$ws.Cells.Item(388, 20).Value = 18

# Row 389 (new): 2a plateado, 3/Feb/2022, Provincia de Melipilla
$ws.Cells.Item(389, 1).Value = 4
$ws.Cells.Item(389, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(389, 3).Value = "Los Lagos"
$ws.Cells.Item(389, 4).Value = 44595
$ws.Cells.Item(389, 5).Value = 10
$ws.Cells.Item(389, 6).Value = "Fruta"
$ws.Cells.Item(389, 7).Value = 100102
$ws.Cells.Item(389, 8).Value = "Cítricos"
$ws.Cells.Item(389, 9).Value = 100102003
$ws.Cells.Item(389, 10).Value = "Limón"
$ws.Cells.Item(389, 11).Value = "Sin especificar"
$ws.Cells.Item(389, 12).Value = "2a plateado"
$ws.Cells.Item(389, 13).Value = 300
$ws.Cells.Item(389, 14).Value = 21000
$ws.Cells.Item(389, 15).Value = 21000
$ws.Cells.Item(389, 16).Value = 21000
$ws.Cells.Item(389, 17).Value = "$/malla 18 kilos"
$ws.Cells.Item(389, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(389, 19).Value = 1167
$ws.Cells.Item(389, 20).Value = 18
